$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, pushing existing rows 6-9 down to 7-10
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly entry
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(6, 3).Value = "Bíobío"
$ws.Cells.Item(6, 4).Value = 45141
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(6, 6).Value = 100112035
$ws.Cells.Item(6, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(6, 11).Value = 8500
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 8800
$ws.Cells.Item(6, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(6, 16).Value = 587
$ws.Cells.Item(6, 17).Value = 15
$ws.Cells.Item(6, 18).Value = "Hortaliza"
